$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1074, pushing the existing rows (1074-1113) down to (1075-1114).
$ws.Rows.Item(1074).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A1074").Value = 8
$ws.Range("B1074").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1074").Value = "Coquimbo"
$ws.Range("D1074").Value = 45075
$ws.Range("E1074").Value = 4
$ws.Range("F1074").Value = 100112023
$ws.Range("G1074").Value = "Brócoli"
$ws.Range("H1074").Value = "Sin especificar"
$ws.Range("I1074").Value = "Primera"
$ws.Range("J1074").Value = 1800
$ws.Range("K1074").Value = 700
$ws.Range("L1074").Value = 800
$ws.Range("M1074").Value = 750
$ws.Range("N1074").Value = "$/unidad"
$ws.Range("O1074").Value = "Provincia del Elquí"
$ws.Range("P1074").Value = 750
$ws.Range("Q1074").Value = 1
$ws.Range("R1074").Value = "Hortaliza"

# Keep the date format consistent with the rest of column D.
$ws.Range("D1074").NumberFormat = "YYYY-MM-DD HH:MM:SS"
